$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text-number-format for price cells whose new values would otherwise be
# auto-coerced into numeric cells by Excel (single-dot decimal strings).
$textForceCells = @("D5", "D9", "D11", "D15", "D18", "D20", "D22", "D25", "D27", "D32", "D36", "D39", "D40", "D42", "D43", "D45", "D48", "D49", "D50", "D51")
foreach ($cellRef in $textForceCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

$ws.Range("D2").Value = "26.946.36"
$ws.Range("E2").Value = "  +0.38%  "
$ws.Range("D3").Value = "1.557.16"
$ws.Range("E3").Value = "  -0.18%  "
$ws.Range("E4").Value = "  +0.43%  "
$ws.Range("D5").Value = "206.86"
$ws.Range("E5").Value = "  +0.74%  "
$ws.Range("E6").Value = "  +1.38%  "
$ws.Range("E7").Value = "  +0.42%  "
$ws.Range("E8").Value = "  +0.32%  "
$ws.Range("D9").Value = "21.51"
$ws.Range("E9").Value = "  -0.26%  "
$ws.Range("E10").Value = "  -0.27%  "
$ws.Range("D11").Value = "0.0860"
$ws.Range("E11").Value = "  -0.12%  "
$ws.Range("D12").Value = "1.777.93"
$ws.Range("E12").Value = "  -0.22%  "
$ws.Range("D13").Value = "1.557.90"
$ws.Range("E13").Value = "  -0.05%  "
$ws.Range("E14").Value = "  -0.26%  "
$ws.Range("D15").Value = "0.515"
$ws.Range("E15").Value = "  +0.39%  "
$ws.Range("D16").Value = "26.941.15"
$ws.Range("E16").Value = "  +0.37%  "
$ws.Range("E17").Value = "  +0.72%  "
$ws.Range("D18").Value = "214.75"
$ws.Range("E18").Value = "  +0.17%  "
$ws.Range("E19").Value = "  +0.54%  "
$ws.Range("D20").Value = "7.26"
$ws.Range("E20").Value = "  -0.81%  "
$ws.Range("E21").Value = "  +0.44%  "
$ws.Range("D22").Value = "4.05"
$ws.Range("E22").Value = "  -1.75%  "
$ws.Range("E23").Value = "  +0.67%  "
$ws.Range("E24").Value = "  -2.48%  "
$ws.Range("D25").Value = "153.54"
$ws.Range("E25").Value = "  +0.24%  "
$ws.Range("E26").Value = "  +1.34%  "
$ws.Range("D27").Value = "14.90"
$ws.Range("E27").Value = "  -0.30%  "
$ws.Range("E28").Value = "  +0.41%  "
$ws.Range("E29").Value = "  +0.84%  "
$ws.Range("E30").Value = "  -1.32%  "
$ws.Range("E31").Value = "  -0.69%  "
$ws.Range("D32").Value = "3.24"
$ws.Range("E32").Value = "  +1.35%  "
$ws.Range("D33").Value = "1.369.98"
$ws.Range("E33").Value = "  -0.42%  "
$ws.Range("E34").Value = "  +1.41%  "
$ws.Range("E35").Value = "  +2.45%  "
$ws.Range("D36").Value = "0.972"
$ws.Range("E36").Value = "  +5.30%  "
$ws.Range("E37").Value = "  +0.51%  "
$ws.Range("E38").Value = "  +0.58%  "
$ws.Range("D39").Value = "0.520"
$ws.Range("E39").Value = "  -1.07%  "
$ws.Range("D40").Value = "0.808"
$ws.Range("E40").Value = "  -0.03%  "
$ws.Range("E41").Value = "  +0.41%  "
$ws.Range("B42").Value = "FraxShare"
$ws.Range("C42").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D42").Value = "5.51"
$ws.Range("E42").Value = "  -0.93%  "
$ws.Range("B43").Value = "WEMIXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D43").Value = "0.981"
$ws.Range("E43").Value = "  -0.84%  "
$ws.Range("E44").Value = "  +3.15%  "
$ws.Range("D45").Value = "63.74"
$ws.Range("E45").Value = "  +0.39%  "
$ws.Range("E46").Value = "  -2.32%  "
$ws.Range("D47").Value = "1.691.37"
$ws.Range("E47").Value = "  -0.29%  "
$ws.Range("D48").Value = "86.04"
$ws.Range("E48").Value = "  -0.49%  "
$ws.Range("D49").Value = "0.0508"
$ws.Range("E49").Value = "  -0.54%  "
$ws.Range("B50").Value = "Algorand"
$ws.Range("C50").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D50").Value = "0.0955"
$ws.Range("E50").Value = "  +0.16%  "
$ws.Range("B51").Value = "USDD"
$ws.Range("C51").Value = "https://coinranking.com/coin/z2PZIKQL7+usdd-usdd"
$ws.Range("D51").Value = "1.01"
$ws.Range("E51").Value = "  +0.46%  "
